$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 111476779
$ws.Range("B32").Value = 89419
$ws.Range("D32").Value = 'NT'
$ws.Range("E32").Value = 1204
$ws.Range("F32").Value = 'Gränsticka'
$ws.Range("G32").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H32").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P32").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q32").Value = 519682.0588881435
$ws.Range("R32").Value = 7151062.961124314
$ws.Range("S32").Value = 10
$ws.Range("AW32").Value = 'Signe Propst'
$ws.Range("AX32").Value = 'Signe Propst'

# Row 33
$ws.Range("A33").Value = 111478332
$ws.Range("B33").Value = 96368
$ws.Range("D33").Value = 'LC'
$ws.Range("E33").Value = 221952
$ws.Range("F33").Value = 'Spindelblomster'
$ws.Range("G33").Value = 'Neottia cordata'
$ws.Range("H33").Value = '(L.) Rich.'
$ws.Range("P33").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q33").Value = 519691.5529315504
$ws.Range("R33").Value = 7150796.076782022
$ws.Range("S33").Value = 10
$ws.Range("AW33").Value = 'Elicia Olsson'
$ws.Range("AX33").Value = 'Elicia Olsson, Astrid Blomberg, Elias Blad, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 34
$ws.Range("A34").Value = 111478318
$ws.Range("B34").Value = 96368
$ws.Range("D34").Value = 'LC'
$ws.Range("E34").Value = 221952
$ws.Range("F34").Value = 'Spindelblomster'
$ws.Range("G34").Value = 'Neottia cordata'
$ws.Range("H34").Value = '(L.) Rich.'
$ws.Range("P34").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q34").Value = 519477.1715801154
$ws.Range("R34").Value = 7151301.40310728
$ws.Range("S34").Value = 10
$ws.Range("AW34").Value = 'Elicia Olsson'
$ws.Range("AX34").Value = 'Elicia Olsson, Astrid Blomberg, Elias Blad, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 35
$ws.Range("A35").Value = 111477350
$ws.Range("B35").Value = 78605
$ws.Range("D35").Value = 'LC'
$ws.Range("E35").Value = 6462
$ws.Range("F35").Value = 'Stuplav'
$ws.Range("G35").Value = 'Nephroma bellum'
$ws.Range("H35").Value = '(Spreng.) Tuck.'
$ws.Range("P35").Value = 'O om Hyktabergets naturreservat, Jmt'
$ws.Range("Q35").Value = 519486.9696307178
$ws.Range("R35").Value = 7151254.464493743
$ws.Range("S35").Value = 15
$ws.Range("AW35").Value = 'Elvira Klang'
$ws.Range("AX35").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 36
$ws.Range("A36").Value = 111477372
$ws.Range("B36").Value = 89686
$ws.Range("D36").Value = 'NT'
$ws.Range("E36").Value = 658
$ws.Range("F36").Value = 'Rosenticka'
$ws.Range("G36").Value = 'Rhodofomes roseus'
$ws.Range("H36").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("P36").Value = 'O om Hyktabergets naturreservat, Jmt'
$ws.Range("Q36").Value = 519504.8857270729
$ws.Range("R36").Value = 7151558.564654102
$ws.Range("S36").Value = 15
$ws.Range("AW36").Value = 'Elvira Klang'
$ws.Range("AX36").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 37
$ws.Range("A37").Value = 111477366
$ws.Range("B37").Value = 78579
$ws.Range("D37").Value = 'NT'
$ws.Range("E37").Value = 2081
$ws.Range("F37").Value = 'Skrovellav'
$ws.Range("G37").Value = 'Lobaria scrobiculata'
$ws.Range("H37").Value = '(Scop.) DC.'
$ws.Range("P37").Value = 'O om Hyktabergets naturreservat, Jmt'
$ws.Range("Q37").Value = 519515.8996614926
$ws.Range("R37").Value = 7151321.485724327
$ws.Range("S37").Value = 15
$ws.Range("AW37").Value = 'Elvira Klang'
$ws.Range("AX37").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 88
$ws.Range("A88").Value = 111477396
$ws.Range("B88").Value = 77515
$ws.Range("D88").Value = 'NT'
$ws.Range("E88").Value = 6425
$ws.Range("F88").Value = 'Garnlav'
$ws.Range("G88").Value = 'Alectoria sarmentosa'
$ws.Range("H88").Value = '(Ach.) Ach.'
$ws.Range("P88").Value = 'O om Hyktabergets naturreservat, Jmt'
$ws.Range("Q88").Value = 519500.1483663829
$ws.Range("R88").Value = 7151151.488534225
$ws.Range("S88").Value = 15
$ws.Range("AW88").Value = 'Elvira Klang'
$ws.Range("AX88").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 89
$ws.Range("A89").Value = 111477364
$ws.Range("B89").Value = 78579
$ws.Range("D89").Value = 'NT'
$ws.Range("E89").Value = 2081
$ws.Range("F89").Value = 'Skrovellav'
$ws.Range("G89").Value = 'Lobaria scrobiculata'
$ws.Range("H89").Value = '(Scop.) DC.'
$ws.Range("P89").Value = 'O om Hyktabergets naturreservat, Jmt'
$ws.Range("Q89").Value = 519703.7324451482
$ws.Range("R89").Value = 7150849.201368837
$ws.Range("S89").Value = 15
$ws.Range("AW89").Value = 'Elvira Klang'
$ws.Range("AX89").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 147
$ws.Range("A147").Value = 111480405
$ws.Range("B147").Value = 89419
$ws.Range("D147").Value = 'NT'
$ws.Range("E147").Value = 1204
$ws.Range("F147").Value = 'Gränsticka'
$ws.Range("G147").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H147").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P147").Value = 'O om Hyktabergets NR, Jmt'
$ws.Range("Q147").Value = 519620.3788537583
$ws.Range("R147").Value = 7151053.076734134
$ws.Range("S147").Value = 15
$ws.Range("AW147").Value = 'Tore Dahlberg'
$ws.Range("AX147").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'

# Row 148
$ws.Range("A148").Value = 111480454
$ws.Range("B148").Value = 77515
$ws.Range("D148").Value = 'NT'
$ws.Range("E148").Value = 6425
$ws.Range("F148").Value = 'Garnlav'
$ws.Range("G148").Value = 'Alectoria sarmentosa'
$ws.Range("H148").Value = '(Ach.) Ach.'
$ws.Range("P148").Value = 'O om Hyktabergets NR, Jmt'
$ws.Range("Q148").Value = 519717.413156116
$ws.Range("R148").Value = 7150803.575370145
$ws.Range("S148").Value = 15
$ws.Range("AW148").Value = 'Tore Dahlberg'
$ws.Range("AX148").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'

# Row 149
$ws.Range("A149").Value = 111480402
$ws.Range("B149").Value = 89419
$ws.Range("D149").Value = 'NT'
$ws.Range("E149").Value = 1204
$ws.Range("F149").Value = 'Gränsticka'
$ws.Range("G149").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H149").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P149").Value = 'O om Hyktabergets NR, Jmt'
$ws.Range("Q149").Value = 519476.848658799
$ws.Range("R149").Value = 7151487.242657971
$ws.Range("S149").Value = 15
$ws.Range("AW149").Value = 'Tore Dahlberg'
$ws.Range("AX149").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'

# Row 150
$ws.Range("A150").Value = 111492392
$ws.Range("B150").Value = 89423
$ws.Range("D150").Value = 'NT'
$ws.Range("E150").Value = 5432
$ws.Range("F150").Value = 'Granticka'
$ws.Range("G150").Value = 'Porodaedalea chrysoloma'
$ws.Range("H150").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("P150").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q150").Value = 519636.82785619
$ws.Range("R150").Value = 7150845.319999835
$ws.Range("S150").Value = 15
$ws.Range("AW150").Value = 'Karl Soler Kinnerbäck'
$ws.Range("AX150").Value = 'Karl Soler Kinnerbäck, Elicia Olsson, Signe Propst, Tore Dahlberg, Melvin Lewin, Elvira Klang, Filippa Paperin, Elias Blad, Iris Elmér, Ivar Anderberg, Kai Strömberg, Astrid Blomberg'

# Row 151
$ws.Range("A151").Value = 111492380
$ws.Range("B151").Value = 89590
$ws.Range("D151").Value = 'VU'
$ws.Range("E151").Value = 48
$ws.Range("F151").Value = 'Lappticka'
$ws.Range("G151").Value = 'Amylocystis lapponica'
$ws.Range("H151").Value = '(Romell) Singer'
$ws.Range("P151").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q151").Value = 519623.8298980477
$ws.Range("R151").Value = 7150852.136615888
$ws.Range("S151").Value = 15
$ws.Range("AW151").Value = 'Karl Soler Kinnerbäck'
$ws.Range("AX151").Value = 'Karl Soler Kinnerbäck, Elicia Olsson, Signe Propst, Tore Dahlberg, Melvin Lewin, Elvira Klang, Filippa Paperin, Elias Blad, Iris Elmér, Ivar Anderberg, Kai Strömberg, Astrid Blomberg'

# Row 152
$ws.Range("A152").Value = 111480258
$ws.Range("B152").Value = 89686
$ws.Range("D152").Value = 'NT'
$ws.Range("E152").Value = 658
$ws.Range("F152").Value = 'Rosenticka'
$ws.Range("G152").Value = 'Rhodofomes roseus'
$ws.Range("H152").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("P152").Value = 'O om Hyktabergets NR, Jmt'
$ws.Range("Q152").Value = 519695.8962852532
$ws.Range("R152").Value = 7150859.069846104
$ws.Range("S152").Value = 15
$ws.Range("AW152").Value = 'Tore Dahlberg'
$ws.Range("AX152").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'

# Row 169
$ws.Range("A169").Value = 111492387
$ws.Range("B169").Value = 89419
$ws.Range("D169").Value = 'NT'
$ws.Range("E169").Value = 1204
$ws.Range("F169").Value = 'Gränsticka'
$ws.Range("G169").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H169").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("P169").Value = 'Öster om Hyktabergets naturreservat, Jmt'
$ws.Range("Q169").Value = 519599.5474416229
$ws.Range("R169").Value = 7150867.937278651
$ws.Range("S169").Value = 15
$ws.Range("AW169").Value = 'Karl Soler Kinnerbäck'
$ws.Range("AX169").Value = 'Karl Soler Kinnerbäck, Elicia Olsson, Signe Propst, Tore Dahlberg, Melvin Lewin, Elvira Klang, Filippa Paperin, Elias Blad, Iris Elmér, Ivar Anderberg, Kai Strömberg, Astrid Blomberg'

# Row 170
$ws.Range("A170").Value = 111480430
$ws.Range("B170").Value = 89423
$ws.Range("D170").Value = 'NT'
$ws.Range("E170").Value = 5432
$ws.Range("F170").Value = 'Granticka'
$ws.Range("G170").Value = 'Porodaedalea chrysoloma'
$ws.Range("H170").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("P170").Value = 'O om Hyktabergets NR, Jmt'
$ws.Range("Q170").Value = 519582.6775320586
$ws.Range("R170").Value = 7150939.848791949
$ws.Range("S170").Value = 15
$ws.Range("AW170").Value = 'Tore Dahlberg'
$ws.Range("AX170").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'
